# Generate Report for Handback
#
# Applies the "handback" localization-status refresh:
#   - Status cells flip from "In Translation" to "Handed back: in sync with en-US"
#     on the Overview sheet and on each language sheet.
#   - Each language sheet's rows gain the "Latest Target File" (J) and
#     "Latest Handback File" (K) values (with J hyperlinked back to the
#     source .md file, same as column A), and the "Latest Handback
#     DateTime" (L) is stamped with the handback timestamp.
#   - The Status / Latest Target File / Latest Handback File columns are
#     widened to comfortably fit the new long filenames.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# Column-width helper: this host quantizes stored column width to the
# nearest 1/6 of a character after a fixed +5/6 offset, so feed it the
# character width whose quantized result lands on the desired value.
function Set-ColWidth($ws, $colIndex, $desiredStoredWidth) {
    $ws.Columns.Item($colIndex).ColumnWidth = $desiredStoredWidth - 0.8333333333333334
}

# ----------------------------------------------------------------------
# Overview sheet: widen the zh-cn / de-de status columns (E, F) and flip
# the status text shown for both rows.
# ----------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

Set-ColWidth $overview 5 29.9777047293527
Set-ColWidth $overview 6 29.9777047293527

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# ----------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de): same shape of edit, different
# handback file names / timestamp.
# ----------------------------------------------------------------------
$languages = @(
    @{ Sheet = "zh-cn"; HandbackSuffix = "zh-cn"; HandbackStamp = "2017-01-03 08:51:44" },
    @{ Sheet = "de-de"; HandbackSuffix = "de-de"; HandbackStamp = "2017-01-03 08:51:56" }
)

$rows = @(
    @{ Row = 2; Guid = "1f8561ec-e3af-4ce1-ab2f-bc7d7a5feaf3"; Hash = "d667ebb0f26b857e533f589dfdcabb1a7b2570c2" },
    @{ Row = 3; Guid = "b7f56c62-5023-44c4-8c61-bb90899180f5"; Hash = "0acae1050ad6e4e5f4370793f54de5e5d8878b15" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    # Widen Status (C), Latest Target File (J) and Latest Handback File (K).
    Set-ColWidth $ws 3 29.9777047293527
    Set-ColWidth $ws 10 40
    Set-ColWidth $ws 11 40

    foreach ($r in $rows) {
        $rowNum = $r.Row
        $mdName = "$($r.Guid).md"
        $targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/e899886868cffddbd4016a41d72e08a7a9dd6a73/e2e/$mdName"
        $handbackFile = "$($r.Guid).$($r.Hash).$($lang.HandbackSuffix).xlf"

        # Status -> handed back.
        $ws.Range("C$rowNum").Value = $newStatus

        # Latest Target File: same file + hyperlink as column A.
        $ws.Range("J$rowNum").Value = $mdName
        $ws.Hyperlinks.Add($ws.Range("J$rowNum"), $targetUrl, "", "", $mdName)

        # Latest Handback File.
        $ws.Range("K$rowNum").Value = $handbackFile

        # Latest Handback DateTime.
        $ws.Range("L$rowNum").Value = $lang.HandbackStamp
    }
}
